$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "a"/"b" columns (F,G) and retarget D1/E1 to
#     SVSPP/Lm (old "BTS #" / "50 % Size at Maturity (cm)" headers are
#     replaced). Order matters for shared-string allocation: "a" and "b"
#     must be interned before "SVSPP"/"Lm" so the unique-string table ends
#     up in the same order as the target workbook.
$ws.Range("F1").Value2 = "a"
$ws.Range("F1").Style = "Normal"
$ws.Range("G1").Value2 = "b"
$ws.Range("G1").Style = "Normal"
$ws.Range("D1").Value2 = "SVSPP"
$ws.Range("E1").Value2 = "Lm"

# --- New juvenile/adult (a/b) coefficients for the two species that have
#     them (Atlantic cod row 3, haddock row 4).
$ws.Range("F3").Value2 = 0.0069
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value2 = 3.08
$ws.Range("G3").Style = "Normal"

$ws.Range("F4").Value2 = 0.0059
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value2 = 3.13
$ws.Range("G4").Style = "Normal"

# --- Materialize the (otherwise blank) F/G cells on the subtotal row (15)
#     and the trailing formatting-only row (19) so the used range/row spans
#     pick up the two new columns there too, without inheriting the F:N
#     column style.
$ws.Range("F15").Value2 = 1
$ws.Range("F15").Value2 = ""
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value2 = 1
$ws.Range("G15").Value2 = ""
$ws.Range("G15").Style = "Normal"

$ws.Range("F19").Value2 = 1
$ws.Range("F19").Value2 = ""
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value2 = 1
$ws.Range("G19").Value2 = ""
$ws.Range("G19").Style = "Normal"

# --- The F:G columns now carry real header/data, so they should drop out
#     of the "style=3, default width" block that used to span F:N -- that
#     block now only covers H:N. ClearFormats on the whole column is what's
#     needed to actually split the <cols> run, but it also stamps a blank
#     cell into every other row of the used range in F:G -- strip those
#     back out afterwards, keeping only the rows that should carry data.
$ws.Columns("F:G").ClearFormats()
$ws.Range("F2:G2").Clear()
$ws.Range("F5:G14").Clear()
$ws.Range("F16:G18").Clear()

# --- Selection moves from the old E23 to H9 (reflecting the new a/b
#     columns the user was working in).
$ws.Range("H9").Select()
